$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: MMSeqs2 -> MMSeqs2_100 (new values)
$ws.Range("A7").Value = "MMSeqs2_100"
$ws.Range("B7").Value = 0.55
$ws.Range("C7").Value = 0.84
$ws.Range("D7").Value = 0.7
$ws.Range("E7").Value = 0.9399999999999999
$ws.Range("F7").Value = 0.55

# Row 8: Mothur -> MMSeqs2_97 (values from old MMSeqs2 row)
$ws.Range("A8").Value = "MMSeqs2_97"
$ws.Range("B8").Value = 0.62
$ws.Range("C8").Value = 0.85
$ws.Range("D8").Value = 0.77
$ws.Range("E8").Value = 0.9
$ws.Range("F8").Value = 0.68

# Row 9: NBC -> Mothur (values from old Mothur row)
$ws.Range("A9").Value = "Mothur"
$ws.Range("B9").Value = 0.44
$ws.Range("C9").Value = 0.64
$ws.Range("D9").Value = 0.59
$ws.Range("E9").Value = 0.68
$ws.Range("F9").Value = 0.52

# Row 10: Qiime2 -> NBC (values updated)
$ws.Range("A10").Value = "NBC"
$ws.Range("B10").Value = 0.53
$ws.Range("C10").Value = 0.73
$ws.Range("D10").Value = 0.6899999999999999
$ws.Range("E10").Value = 0.78
$ws.Range("F10").Value = 0.6

# Row 11: TNT -> Qiime2 (values from old Qiime2 row)
$ws.Range("A11").Value = "Qiime2"
$ws.Range("B11").Value = 0.32
$ws.Range("C11").Value = 0.62
$ws.Range("D11").Value = 0.48
$ws.Range("E11").Value = 0.6899999999999999
$ws.Range("F11").Value = 0.37

# Row 12: VSEARCH stays, but Accuracy changes
$ws.Range("B12").Value = 0.43
